$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 7 (10006 / "Berühmter Streamer..." entry)
$ws.Rows.Item(7).Delete()

# Delete rows 28 and 29 (10031 / "Gut gemacht..." and 10032 / "Volles Haus..." entries)
# After deleting original row 7, these shift up by one: original 28/29 -> now 27/28
$ws.Rows.Item(27).Delete()
$ws.Rows.Item(27).Delete()

# Update selection to A24
$ws.Range("A24").Select()
